# Weekly data refresh: a new price-report row for Feria Lagunitas de Puerto
# Montt - Cebollín is inserted at row 40, pushing every subsequent record
# down by one row (the last record lands on the newly created row 174).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40:173 down to 41:174, leaving a blank row 40 for the new entry.
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with this week's report.
$ws.Cells.Item(40, 1).Value = 4
$ws.Cells.Item(40, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(40, 3).Value = "Los Lagos"
$ws.Cells.Item(40, 4).Value = 44487
$ws.Cells.Item(40, 5).Value = 10
$ws.Cells.Item(40, 6).Value = 100112037
$ws.Cells.Item(40, 7).Value = "Cebollín"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 60
$ws.Cells.Item(40, 11).Value = 5500
$ws.Cells.Item(40, 12).Value = 6000
$ws.Cells.Item(40, 13).Value = 5750
$ws.Cells.Item(40, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(40, 15).Value = "Región Metropolitana"
$ws.Cells.Item(40, 16).Value = 160
$ws.Cells.Item(40, 17).Value = 36
$ws.Cells.Item(40, 18).Value = "Hortaliza"
